# Update cryptocurrency Price (D) and Volume(1h) (E) columns per the latest scrape.
# A leading apostrophe forces Excel to store the value as text (matching the
# existing inline-string cell type) instead of auto-converting number-looking
# strings like "562.29" into numeric cells; resetting Style to "Normal"
# afterwards clears the quote-prefix formatting Excel applies when it sees the
# leading apostrophe, so the cell's style stays identical to before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.712.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.46%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.402.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.87%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'562.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.62%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'176.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.96%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.393.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.90%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.12%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.173"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.63%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.00%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'54.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.87%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000281"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.72%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'9.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.12%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.938.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.68%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'18.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.23%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +1.19%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.393.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.77%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'65.608.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.37%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.06%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.40%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'463.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.81%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.06%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'4.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.30%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'14.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +5.96%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'87.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.68%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.87%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.55%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.39%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'31.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.35%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.52%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'63.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +6.81%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.08%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'582.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.10%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -0.16%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.09%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +4.14%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.76%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'36.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.28%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.376"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0₃0746"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.85%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.108.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.06%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +1.43%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.37%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -2.92%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +1.89%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.07%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.11%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'140.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.57%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.19%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'8.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.14%  "
$ws.Range("E51").Style = "Normal"
